# Revised implementation of link process partitioning and added DSI to air pollution controls
#
# The apcd_abbrev lookup table maps an "original_eia_860" abbreviation (column C)
# to the new "apcd_equip" abbreviation (columns A/B). Row 19 (DSI / "Dry sorbent
# (powder) injection type (DSI)") was previously partitioned/linked to the
# "dFGD" process group; it is revised here to link directly to "DSI" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the linked/partitioned process abbreviation for the DSI row.
$ws.Range("C19").Value = "DSI"

# Reflect the resulting cursor/viewport position left behind by the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C20").Select()
